$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4384205043315887
$ws.Range("B1").Value = 0.8841232657432556
$ws.Range("C1").Value = 4.397961616516113
$ws.Range("D1").Value = 2.194677591323853
$ws.Range("E1").Value = 0.7833092212677002
